$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.756.51'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +0.19%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.636.63'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.32%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '216.67'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -1.11%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.507'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +1.99%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.57%  '
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +2.15%  '
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +0.62%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.81'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +4.88%  '
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +0.11%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.865.80'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +0.37%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.619.90'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -0.10%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.11'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +0.25%  '
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +1.65%  '
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +3.33%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '26.762.00'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +0.31%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.0₃0727'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +0.02%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '217.86'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +2.62%  '
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -0.69%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.65'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +7.00%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.38'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +1.41%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.43'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +4.34%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.13'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +1.42%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '147.16'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -0.43%  '
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -0.76%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.38'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +5.13%  '
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +0.79%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.72'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +1.14%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0504'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +0.68%  '
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -1.32%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.32'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -1.46%  '
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +1.40%  '
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +1.90%  '
$ws.Range("B35").NumberFormat = "@"
$ws.Range("B35").Value = 'Maker'
$ws.Range("C35").NumberFormat = "@"
$ws.Range("C35").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.254.42'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +0.11%  '
$ws.Range("B36").NumberFormat = "@"
$ws.Range("B36").Value = 'HuobiToken'
$ws.Range("C36").NumberFormat = "@"
$ws.Range("C36").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.44'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -0.25%  '
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +0.87%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.532'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +1.96%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.832'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +3.70%  '
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -0.62%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.43'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +3.44%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.777.56'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '61.65'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +3.64%  '
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -0.86%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '91.55'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -0.33%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.56'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -0.25%  '
$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = 'Cronos'
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0513'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -0.53%  '
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.61'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +1.95%  '
$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = 'Algorand'
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0964'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +1.00%  '
$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = 'USDD'
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.00'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -0.63%  '
